{"js": "// Replace the 25 \"A\u00d7B=C\" arithmetic-answer strings in the table cells with\n// their updated values, per the commit diff. Each old value is unique in the\n// document, so an exact (non-wildcard, case-sensitive, whole-match) search-\n// and-replace for each pair is safe and unambiguous.\nconst replacements = [\n  [\"522\u00d73=1566\", \"201\u00d76=1206\"],\n  [\"319\u00d79=2871\", \"937\u00d76=5622\"],\n  [\"514\u00d73=1542\", \"274\u00d77=1918\"],\n  [\"967\u00d75=4835\", \"250\u00d78=2000\"],\n  [\"766\u00d72=1532\", \"947\u00d78=7576\"],\n  [\"421\u00d79=3789\", \"718\u00d78=5744\"],\n  [\"541\u00d73=1623\", \"448\u00d76=2688\"],\n  [\"242\u00d76=1452\", \"625\u00d78=5000\"],\n  [\"466\u00d74=1864\", \"713\u00d72=1426\"],\n  [\"105\u00d74=420\", \"973\u00d75=4865\"],\n  [\"823\u00d72=1646\", \"113\u00d77=791\"],\n  [\"537\u00d72=1074\", \"984\u00d76=5904\"],\n  [\"782\u00d75=3910\", \"349\u00d75=1745\"],\n  [\"706\u00d77=4942\", \"535\u00d72=1070\"],\n  [\"910\u00d78=7280\", \"803\u00d75=4015\"],\n  [\"632\u00d73=1896\", \"360\u00d79=3240\"],\n  [\"860\u00d77=6020\", \"589\u00d72=1178\"],\n  [\"476\u00d72=952\", \"647\u00d78=5176\"],\n  [\"939\u00d79=8451\", \"920\u00d73=2760\"],\n  [\"484\u00d79=4356\", \"633\u00d75=3165\"],\n  [\"142\u00d78=1136\", \"291\u00d79=2619\"],\n  [\"806\u00d79=7254\", \"418\u00d73=1254\"],\n  [\"970\u00d75=4850\", \"624\u00d73=1872\"],\n  [\"500\u00d74=2000\", \"214\u00d73=642\"],\n  [\"624\u00d78=4992\", \"786\u00d79=7074\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"A\u00d7B=C\" arithmetic-answer strings in the table cells with\n# their updated values, per the commit diff. Each old value is unique in the\n# document, so a plain (non-wildcard, case-sensitive) Find/Replace per pair\n# is safe and unambiguous.\n$replacements = @(\n  @(\"522\u00d73=1566\", \"201\u00d76=1206\"),\n  @(\"319\u00d79=2871\", \"937\u00d76=5622\"),\n  @(\"514\u00d73=1542\", \"274\u00d77=1918\"),\n  @(\"967\u00d75=4835\", \"250\u00d78=2000\"),\n  @(\"766\u00d72=1532\", \"947\u00d78=7576\"),\n  @(\"421\u00d79=3789\", \"718\u00d78=5744\"),\n  @(\"541\u00d73=1623\", \"448\u00d76=2688\"),\n  @(\"242\u00d76=1452\", \"625\u00d78=5000\"),\n  @(\"466\u00d74=1864\", \"713\u00d72=1426\"),\n  @(\"105\u00d74=420\", \"973\u00d75=4865\"),\n  @(\"823\u00d72=1646\", \"113\u00d77=791\"),\n  @(\"537\u00d72=1074\", \"984\u00d76=5904\"),\n  @(\"782\u00d75=3910\", \"349\u00d75=1745\"),\n  @(\"706\u00d77=4942\", \"535\u00d72=1070\"),\n  @(\"910\u00d78=7280\", \"803\u00d75=4015\"),\n  @(\"632\u00d73=1896\", \"360\u00d79=3240\"),\n  @(\"860\u00d77=6020\", \"589\u00d72=1178\"),\n  @(\"476\u00d72=952\", \"647\u00d78=5176\"),\n  @(\"939\u00d79=8451\", \"920\u00d73=2760\"),\n  @(\"484\u00d79=4356\", \"633\u00d75=3165\"),\n  @(\"142\u00d78=1136\", \"291\u00d79=2619\"),\n  @(\"806\u00d79=7254\", \"418\u00d73=1254\"),\n  @(\"970\u00d75=4850\", \"624\u00d73=1872\"),\n  @(\"500\u00d74=2000\", \"214\u00d73=642\"),\n  @(\"624\u00d78=4992\", \"786\u00d79=7074\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Text = $oldText\n  $rng.Find.Replacement.Text = $newText\n  $rng.Find.Forward = $true\n  $rng.Find.Wrap = 1\n  $rng.Find.Format = $false\n  $rng.Find.MatchCase = $true\n  $rng.Find.MatchWholeWord = $false\n  $rng.Find.MatchWildcards = $false\n  $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
